# Apply the edits described by the commit:
# "fix the grammar of dataset section, and minor mistakes pointed by reviewer 1"
#
# This touches the "edge_for_each_paper_accession" worksheet:
#  - removes a block of blank rows (20-23, 25) so the trailing rows shift up
#  - adds two new corrected-correlation cells (H8, H10) with formulas
#  - narrows column B and adds an explicit width for the new column H
#  - updates the view (zoom / scroll / selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("edge_for_each_paper_accession")
$ws.Activate()

# --- remove the blank rows in between, shifting rows 24/26/27 up to 20/21/22 ---
$ws.Rows("25:25").Delete()
$ws.Rows("20:23").Delete()

# --- add the two new "corrected" correlation columns ---
$ws.Range("H8").Formula = "=G8^2*(800-1)/(800-1-1)"
$ws.Range("H10").Formula = "=D14^2*(800-1)/(800-1-3)"

# --- column widths ---
$ws.Columns("B").ColumnWidth = 14.7714285714286
$ws.Columns("H").ColumnWidth = 12.8571428571429

# --- sheet view: zoom, scroll position, and selection ---
$win = $excel.ActiveWindow
$win.Zoom = 175
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("F19").Select()
